# Auto-generated edit script applying crypto price/volume updates
# to match the commit "Updated cryptos list on Sat Dec 23 18:15:55 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force the cell to be stored as text (matching the original inlineStr cells),
    # then strip the temporary Text number-format so no stray style is introduced.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextCell $ws.Range("D2") '43.815.39'
Set-TextCell $ws.Range("E2") '  -0.02%  '
Set-TextCell $ws.Range("D3") '2.291.74'
Set-TextCell $ws.Range("E3") '  -1.97%  '
Set-TextCell $ws.Range("E4") '  +0.07%  '
Set-TextCell $ws.Range("D5") '98.75'
Set-TextCell $ws.Range("E5") '  +2.21%  '
Set-TextCell $ws.Range("D6") '270.68'
Set-TextCell $ws.Range("E6") '  -0.27%  '
Set-TextCell $ws.Range("D7") '0.625'
Set-TextCell $ws.Range("E7") '  -0.27%  '
Set-TextCell $ws.Range("E8") '  +0.08%  '
Set-TextCell $ws.Range("D9") '0.605'
Set-TextCell $ws.Range("E9") '  -3.16%  '
Set-TextCell $ws.Range("D10") '45.04'
Set-TextCell $ws.Range("E10") '  -1.33%  '
Set-TextCell $ws.Range("D11") '0.0932'
Set-TextCell $ws.Range("E11") '  -1.65%  '
Set-TextCell $ws.Range("D12") '7.90'
Set-TextCell $ws.Range("E12") '  -2.60%  '
Set-TextCell $ws.Range("E13") '  +1.55%  '
Set-TextCell $ws.Range("D14") '15.77'
Set-TextCell $ws.Range("E14") '  +0.79%  '
Set-TextCell $ws.Range("D15") '2.635.81'
Set-TextCell $ws.Range("E15") '  -1.79%  '
Set-TextCell $ws.Range("D16") '0.852'
Set-TextCell $ws.Range("E16") '  -1.93%  '
Set-TextCell $ws.Range("D17") '2.290.35'
Set-TextCell $ws.Range("E17") '  -1.96%  '
Set-TextCell $ws.Range("D18") '43.751.74'
Set-TextCell $ws.Range("E18") '  -0.01%  '
Set-TextCell $ws.Range("D19") '0.0000110'
Set-TextCell $ws.Range("E19") '  +0.56%  '
Set-TextCell $ws.Range("D20") '6.19'
Set-TextCell $ws.Range("E20") '  -4.35%  '
Set-TextCell $ws.Range("D21") '72.10'
Set-TextCell $ws.Range("E21") '  -0.86%  '
Set-TextCell $ws.Range("D22") '2.47'
Set-TextCell $ws.Range("E22") '  +8.35%  '
Set-TextCell $ws.Range("D23") '232.98'
Set-TextCell $ws.Range("E23") '  -2.66%  '
Set-TextCell $ws.Range("D24") '2.84'
Set-TextCell $ws.Range("E24") '  +12.06%  '
Set-TextCell $ws.Range("E25") '  -3.43%  '
Set-TextCell $ws.Range("E26") '  -0.03%  '
Set-TextCell $ws.Range("D27") '11.21'
Set-TextCell $ws.Range("E27") '  -1.84%  '
Set-TextCell $ws.Range("E28") '  -0.95%  '
Set-TextCell $ws.Range("D29") '2.29'
Set-TextCell $ws.Range("E29") '  +0.81%  '
Set-TextCell $ws.Range("D30") '38.24'
Set-TextCell $ws.Range("E30") '  +0.61%  '
Set-TextCell $ws.Range("D31") '176.51'
Set-TextCell $ws.Range("E31") '  +1.65%  '
Set-TextCell $ws.Range("D32") '21.79'
Set-TextCell $ws.Range("E32") '  -3.40%  '
Set-TextCell $ws.Range("D33") '0.0890'
Set-TextCell $ws.Range("E33") '  -1.43%  '
Set-TextCell $ws.Range("D34") '5.41'
Set-TextCell $ws.Range("E34") '  -1.67%  '
Set-TextCell $ws.Range("E35") '  +0.41%  '
Set-TextCell $ws.Range("E36") '  +7.50%  '
Set-TextCell $ws.Range("E37") '  -0.44%  '
Set-TextCell $ws.Range("D38") '0.0350'
Set-TextCell $ws.Range("E38") '  -3.07%  '
Set-TextCell $ws.Range("D39") '3.54'
Set-TextCell $ws.Range("E39") '  +4.01%  '
Set-TextCell $ws.Range("E40") '  -0.44%  '
Set-TextCell $ws.Range("E41") '  -2.53%  '
Set-TextCell $ws.Range("E42") '  -1.45%  '
Set-TextCell $ws.Range("D43") '12.17'
Set-TextCell $ws.Range("E43") '  +0.10%  '
Set-TextCell $ws.Range("D44") '64.55'
Set-TextCell $ws.Range("E44") '  +3.27%  '
Set-TextCell $ws.Range("D45") '8.85'
Set-TextCell $ws.Range("E45") '  -4.24%  '
Set-TextCell $ws.Range("E46") '  -3.08%  '
Set-TextCell $ws.Range("E47") '  -1.80%  '
Set-TextCell $ws.Range("E48") '  +0.94%  '
Set-TextCell $ws.Range("D49") '97.86'
Set-TextCell $ws.Range("E49") '  -2.66%  '
Set-TextCell $ws.Range("E50") '  +11.65%  '
Set-TextCell $ws.Range("D51") '0.439'
Set-TextCell $ws.Range("E51") '  +5.48%  '
